# Update "想去人数" (want-to-go count, column F) values on the "展览"
# and "全部类型" sheets to reflect freshly scraped totals.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Row -> new F value, for the 展览 sheet
$exhibitUpdates = @{
    3  = 558
    4  = 1592
    6  = 210
    7  = 784
    8  = 1064
    10 = 378
    12 = 529
    13 = 38
    14 = 6608
    15 = 88
    20 = 1048
    21 = 15823
    22 = 1563
    27 = 11195
    29 = 4392
    30 = 276
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Cells.Item($row, 6).Value = $exhibitUpdates[$row]
}

# Row -> new F value, for the 全部类型 sheet
$allUpdates = @{
    3  = 558
    4  = 1592
    6  = 210
    7  = 784
    9  = 1064
    11 = 378
    13 = 529
    15 = 38
    16 = 6608
    17 = 88
    23 = 1048
    24 = 15823
    25 = 1563
    31 = 11195
    33 = 4392
    34 = 276
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $allUpdates[$row]
}
